$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.023.62"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3
$ws.Range("D3").Value = "2.305.36"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.60%  "

# Row 7
$ws.Range("E7").Value = "  -1.72%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("E9").Value = "  -0.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.58%  "

# Row 13
$ws.Range("E13").Value = "  +1.82%  "

# Row 14
$ws.Range("E14").Value = "  +2.43%  "

# Row 15
$ws.Range("D15").Value = "2.663.90"
$ws.Range("E15").Value = "  +0.60%  "

# Row 16
$ws.Range("D16").Value = "2.327.01"
$ws.Range("E16").Value = "  +1.58%  "

# Row 17
$ws.Range("E17").Value = "  +1.29%  "

# Row 18
$ws.Range("D18").Value = "42.893.24"
$ws.Range("E18").Value = "  +0.38%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("E20").Value = "  -0.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.04%  "

# Row 25
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.58%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.39%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
$ws.Range("E35").Value = "  -7.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0690"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.95%  "

# Row 38
$ws.Range("E38").Value = "  -0.36%  "

# Row 39
$ws.Range("E39").Value = "  +0.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.77%  "

# Row 41
$ws.Range("E41").Value = "  -0.51%  "

# Row 42
$ws.Range("D42").Value = "1.999.74"
$ws.Range("E42").Value = "  -0.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.49%  "

# Row 44
$ws.Range("E44").Value = "  -0.70%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.23%  "

# Row 47
$ws.Range("E47").Value = "  +0.29%  "

# Row 48
$ws.Range("E48").Value = "  -4.14%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.22%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.532.18"
$ws.Range("E50").Value = "  +0.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
